# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last refreshed" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 06:19"

# 2. India (row 5) - refreshed case counts
$ws.Range("B5").Value = 7494551
$ws.Range("C5").Value = 1824
$ws.Range("D5").Value = 6597209
$ws.Range("E5").Value = 783278

# 3. Venezuela (row 55) - refreshed case counts
$ws.Range("B55").Value = 86289
$ws.Range("D55").Value = 78847
$ws.Range("E55").Value = 6711
$ws.Range("H55").Value = 731

# 4. Tailandia (row 147) - refreshed case counts
$ws.Range("B147").Value = 3686
$ws.Range("C147").Value = 7
$ws.Range("D147").Value = 3481
$ws.Range("E147").Value = 146

# 5. Mongolia / Butan (rows 186-187) swapped order (Butan overtook Mongolia in
#    total cases) with refreshed data - row 186 now holds Butan, row 187 Mongolia
$ws.Range("A186").Value = "Butan"
$ws.Range("B186").Value = 325
$ws.Range("C186").Value = 9
$ws.Range("D186").Value = 299
$ws.Range("E186").Value = 26

$ws.Range("A187").Value = "Mongolia"
$ws.Range("B187").Value = 320
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 311
$ws.Range("E187").Value = 9
